$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(9).Insert()

$ws.Range("A9").Value = 1
$ws.Range("B9").Value = 2
$ws.Range("C9").Value = 3
$ws.Range("D9").Value = 4
$ws.Range("E9").Value = 5
$ws.Range("F9").Value = 6
$ws.Range("G9").Value = 7
